{"js": "// Hybrid bold + color highlighting for quantitative impact metrics\n// (percentages, dollar amounts, large numbers) across achievement /\n// responsibility bullet paragraphs.\n//\n// For every target paragraph we locate it by a unique substring of its\n// original text, then run a paragraph-scoped search for each metric\n// token (in left-to-right order) and apply bold + the brand navy color\n// (#2C3E50) to just that token's Range. Word/Office.js automatically\n// splits the paragraph's run(s) around the matched Range, which mirrors\n// the run-splitting pattern seen in the target OOXML diff.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// [uniqueSubstringToLocateParagraph, [metricTokensInOrder]]\nconst targets = [\n  [\n    \"Discovered systematic race coding errors affecting all Black and Asian-American voters\",\n    [\"23%\", \"64%\"],\n  ],\n  [\n    \"Utilized advanced sampling methods to decrease survey margin of error from\",\n    [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"],\n  ],\n  [\n    \"Trigonometric algorithm for boundary estimation reduced mapping costs by\",\n    [\"73.5%\", \"$4.7M\"],\n  ],\n  [\n    \"Built real-time FEC analysis systems using Python, Pandas and PySpark\",\n    [\"$2\"],\n  ],\n  [\n    \"Modernized legacy ETL processes by implementing dbt and PySpark workflows\",\n    [\"57%\"],\n  ],\n  [\n    \"178% accuracy improvement in racial classification algorithms\",\n    [\"178%\"],\n  ],\n  [\n    \"Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs\",\n    [\"73.5%\"],\n  ],\n  [\n    \"$4.7M savings enabled nonprofit access\",\n    [\"$4.7M\"],\n  ],\n  [\n    \"Platform impact: Built redistricting system serving\",\n    [\"12,847\"],\n  ],\n  [\n    \"Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from\",\n    [\"\u00b14.2%\", \"\u00b12.1%\"],\n  ],\n  [\n    \"Increased voter turnout prediction accuracy from\",\n    [\"71%\", \"87%\"],\n  ],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const [needle, tokens] of targets) {\n  const para = paragraphs.items.find((p) => p.text.indexOf(needle) !== -1);\n  if (!para) {\n    throw new Error(\"Could not locate target paragraph for: \" + needle);\n  }\n\n  for (const token of tokens) {\n    const found = para.search(token, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n\n    if (found.items.length === 0) {\n      throw new Error(\"Could not locate token '\" + token + \"' in paragraph: \" + needle);\n    }\n\n    // Metric tokens are unique within their paragraph, so the first hit\n    // is the one we want.\n    const hit = found.items[0];\n    hit.font.bold = true;\n    hit.font.color = HIGHLIGHT_COLOR;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Hybrid bold + color highlighting for quantitative impact metrics\n# (percentages, dollar amounts, large numbers) across achievement /\n# responsibility bullet paragraphs.\n#\n# For every target paragraph we locate it by a unique, literal substring\n# of its original text, then run a paragraph-scoped Find.Execute for each\n# metric token (in left-to-right order) and apply Bold + the brand navy\n# color (#2C3E50 -> wdColor 5258796) to just that token's Range. Word\n# automatically splits the paragraph's run(s) around the matched Range,\n# which mirrors the run-splitting pattern in the target OOXML.\n\n# #2C3E50 packed as a Word \"wdColor\" long (BGR byte order: R + G*256 + B*65536)\n$HighlightColor = 5258796\n$PlusMinus = [char]0xB1\n\nfunction Find-UniqueParagraph($doc, [string]$needle) {\n    $paras = $doc.Paragraphs\n    $n = $paras.Count\n    $hit = $null\n    $count = 0\n    for ($i = 1; $i -le $n; $i++) {\n        $p = $paras.Item($i)\n        if ($p.Range.Text.Contains($needle)) {\n            $count++\n            $hit = $p\n        }\n    }\n    if ($count -ne 1) {\n        throw \"Expected exactly 1 paragraph match for '$needle', found $count\"\n    }\n    return $hit\n}\n\nfunction Set-MetricHighlight($para, [string]$token) {\n    $r = $para.Range\n    $r.Find.ClearFormatting()\n    $r.Find.Forward = $true\n    $r.Find.MatchCase = $true\n    $r.Find.MatchWildcards = $false\n    $found = $r.Find.Execute($token)\n    if (-not $found) {\n        throw \"Could not locate token '$token' in paragraph\"\n    }\n    $r.Font.Bold = 1\n    $r.Font.Color = $HighlightColor\n}\n\n$doc = $word.ActiveDocument\n\n# 1) Discovered systematic race coding errors ... from 23% to 64%\n$p1 = Find-UniqueParagraph $doc \"Discovered systematic race coding errors affecting all Black and Asian-American voters\"\nSet-MetricHighlight $p1 '23%'\nSet-MetricHighlight $p1 '64%'\n\n# 2) Utilized advanced sampling methods ... \u00b14.2% to \u00b12.1% ... 71% to 87%\n$needle2 = [string]::Concat(\"margin of error from \", $PlusMinus, \"4.2% to \", $PlusMinus, \"2.1%, increasing\")\n$p2 = Find-UniqueParagraph $doc $needle2\n$tok2a = [string]::Concat($PlusMinus, '4.2%')\n$tok2b = [string]::Concat($PlusMinus, '2.1%')\nSet-MetricHighlight $p2 $tok2a\nSet-MetricHighlight $p2 $tok2b\nSet-MetricHighlight $p2 '71%'\nSet-MetricHighlight $p2 '87%'\n\n# 3) Trigonometric algorithm ... 73.5% ... $4.7M\n$p3 = Find-UniqueParagraph $doc \"Trigonometric algorithm for boundary estimation reduced mapping costs by\"\nSet-MetricHighlight $p3 '73.5%'\nSet-MetricHighlight $p3 '$4.7M'\n\n# 4) Built real-time FEC analysis systems ... valued over $2 trillion\n$p4 = Find-UniqueParagraph $doc \"Built real-time FEC analysis systems using Python, Pandas and PySpark\"\nSet-MetricHighlight $p4 '$2'\n\n# 5) Modernized legacy ETL processes ... by 57%\n$p5 = Find-UniqueParagraph $doc \"Modernized legacy ETL processes by implementing dbt and PySpark workflows\"\nSet-MetricHighlight $p5 '57%'\n\n# 6) 178% accuracy improvement in racial classification algorithms\n$p6 = Find-UniqueParagraph $doc \"178% accuracy improvement in racial classification algorithms\"\nSet-MetricHighlight $p6 '178%'\n\n# 7) Algorithmic innovation ... mapping costs 73.5%\n$p7 = Find-UniqueParagraph $doc \"Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs\"\nSet-MetricHighlight $p7 '73.5%'\n\n# 8) $4.7M savings enabled nonprofit access\n$p8 = Find-UniqueParagraph $doc \"savings enabled nonprofit access\"\nSet-MetricHighlight $p8 '$4.7M'\n\n# 9) Platform impact: Built redistricting system serving 12,847 analysts ...\n$p9 = Find-UniqueParagraph $doc \"Platform impact: Built redistricting system serving\"\nSet-MetricHighlight $p9 '12,847'\n\n# 10) Predictive excellence: ... \u00b14.2% to \u00b12.1%\n$needle10 = [string]::Concat(\"Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from \", $PlusMinus, \"4.2% to \", $PlusMinus, \"2.1%\")\n$p10 = Find-UniqueParagraph $doc $needle10\n$tok10a = [string]::Concat($PlusMinus, '4.2%')\n$tok10b = [string]::Concat($PlusMinus, '2.1%')\nSet-MetricHighlight $p10 $tok10a\nSet-MetricHighlight $p10 $tok10b\n\n# 11) Increased voter turnout prediction accuracy from 71% to 87%\n$p11 = Find-UniqueParagraph $doc \"Increased voter turnout prediction accuracy from\"\nSet-MetricHighlight $p11 '71%'\nSet-MetricHighlight $p11 '87%'\n"}
